$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.081.21"
$ws.Range("E2").Value = "  +6.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.238.41"
$ws.Range("E3").Value = "  +3.06%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "394.56"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.03"
$ws.Range("E6").Value = "  -2.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.576"
$ws.Range("E7").Value = "  +5.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.234.35"
$ws.Range("E8").Value = "  +3.03%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.620"
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "38.95"
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0957"
$ws.Range("E12").Value = "  +9.84%  "
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.750.32"
$ws.Range("E14").Value = "  +2.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.16"
$ws.Range("E15").Value = "  +1.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.06"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.253.18"
$ws.Range("E17").Value = "  +3.33%  "
$ws.Range("E18").Value = "  -2.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.82"
$ws.Range("E19").Value = "  +3.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "56.912.88"
$ws.Range("E20").Value = "  +6.54%  "
$ws.Range("E21").Value = "  +2.22%  "
$ws.Range("E22").Value = "  +8.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.94"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "296.05"
$ws.Range("E24").Value = "  +9.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.98"
$ws.Range("E25").Value = "  +4.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.16"
$ws.Range("E26").Value = "  -2.28%  "
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.68"
$ws.Range("E28").Value = "  -3.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.30"
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.40"
$ws.Range("E32").Value = "  +4.29%  "
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "38.24"
$ws.Range("E34").Value = "  +3.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0484"
$ws.Range("E35").Value = "  -3.12%  "
$ws.Range("E36").Value = "  +1.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.82"
$ws.Range("E37").Value = "  +2.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.52"
$ws.Range("E38").Value = "  -3.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("E40").Value = "  +5.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "135.28"
$ws.Range("E41").Value = "  +3.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.121"
$ws.Range("E42").Value = "  +3.10%  "
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.94"
$ws.Range("E44").Value = "  -3.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.96"
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.281"
$ws.Range("E46").Value = "  -2.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.09"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.157.35"
$ws.Range("E48").Value = "  +3.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.11"
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.99"
$ws.Range("E50").Value = "  +20.27%  "
$ws.Range("E51").Value = "  -3.94%  "
